# Add 10 new question paragraphs to the end of the Computer Graphics
# Question Bank list, matching the "Added 10 more questions" commit.
#
# Strategy: Word's Range.InsertXML() accepts a full WordOpenXML "flat
# package" document and, when invoked on a Range that falls inside the
# last paragraph of the body, it REPLACES that containing paragraph with
# the supplied <w:p> elements -- but it also always leaves one extra,
# empty trailing paragraph behind (mirroring the tail of the paragraph
# that got split). So: we (1) re-supply the original last paragraph
# unchanged, followed by the 10 new paragraphs, in one InsertXML call,
# then (2) remove the single stray empty paragraph it leaves behind by
# deleting the paragraph-mark that separates our real last paragraph
# from it (that merge is what actually collapses the count back down).

$d = $word.ActiveDocument

$degree = [char]0x00B0

$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/></w:pPr>'

$xmlHead = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlTail = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Original last paragraph, reproduced unchanged.
$origPara = '<w:p>' + $pPr + '<w:r><w:t>A point (4</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>,3</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>) is rotated counter</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>clockwise by an angle of45' + $degree + '. Find the rotation matrix and the resultant point</w:t></w:r></w:p>'

$q1  = '<w:p>' + $pPr + '<w:r><w:t>Name any three font editing tools.</w:t></w:r></w:p>'
$q2  = '<w:p>' + $pPr + '<w:r><w:t>Differentiate serif and sans serif fonts.</w:t></w:r></w:p>'
$q3  = '<w:p>' + $pPr + '<w:r><w:t>Distinguish between window port &amp; view port?</w:t></w:r></w:p>'
$q4  = '<w:p>' + $pPr + '<w:r><w:t>Define clipping?</w:t></w:r></w:p>'
$q5  = '<w:p>' + $pPr + '<w:r><w:t>What is the need of homogeneous coordinates?</w:t></w:r></w:p>'
$q6  = '<w:p>' + $pPr + '<w:r><w:t>Distinguish between uniform scaling and differential</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>scaling?</w:t></w:r></w:p>'
$q7  = '<w:p>' + $pPr + '<w:r><w:t>What is fixed point scaling?</w:t></w:r></w:p>'
$q8  = '<w:p>' + $pPr + '<w:r><w:t>What is Bezier Basis Function?</w:t></w:r></w:p>'
$q9  = '<w:p>' + $pPr + '<w:r><w:t xml:space="preserve">What is surface </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>patch</w:t></w:r><w:r><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>'
$q10 = '<w:p>' + $pPr + '<w:r><w:t>Define B-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Spline</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> curve?</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>'

$fragment = $xmlHead + $origPara + $q1 + $q2 + $q3 + $q4 + $q5 + $q6 + $q7 + $q8 + $q9 + $q10 + $xmlTail

$lastPara = $d.Paragraphs.Last
$target = $lastPara.Range.Duplicate
$target.Collapse(0)
$target.InsertXML($fragment)

# InsertXML leaves one stray empty paragraph after our new content; merge
# it away by deleting the paragraph mark that precedes it.
$n = $d.Paragraphs.Count
$mergePos = $d.Paragraphs.Item($n - 1).Range.End
$markRange = $d.Range($mergePos - 1, $mergePos)
$markRange.Delete()

Write-Host "Final paragraph count:" $d.Paragraphs.Count
